$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings are not
# converted to numbers by Excel (matches original inlineStr typing).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '67.460.96'
$ws.Range("E2").Value = '  +0.33%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '2.476.89'
$ws.Range("E3").Value = '  -2.89%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '588.79'
$ws.Range("E5").Value = '  -0.31%  '

$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '172.66'
$ws.Range("E6").Value = '  -0.53%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = '0.522'
$ws.Range("E8").Value = '  -1.80%  '

$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D9").Value = '2.479.23'
$ws.Range("E9").Value = '  -2.84%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.139'
$ws.Range("E10").Value = '  -0.04%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '0.166'
$ws.Range("E11").Value = '  +1.38%  '

$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D12").Value = '5.08'
$ws.Range("E12").Value = '  -1.44%  '

$ws.Range("B13").Value = 'Cardano'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D13").Value = '0.339'
$ws.Range("E13").Value = '  -2.84%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '26.15'
$ws.Range("E14").Value = '  -3.14%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.955.54'
$ws.Range("E15").Value = '  -1.86%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.0000176'
$ws.Range("E16").Value = '  -0.99%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '67.368.32'
$ws.Range("E17").Value = '  +0.41%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.555.80'
$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '11.69'
$ws.Range("E19").Value = '  +3.05%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '7.96'
$ws.Range("E20").Value = '  -1.12%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '366.04'
$ws.Range("E21").Value = '  +2.50%  '

$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").Value = '4.12'
$ws.Range("E22").Value = '  -2.21%  '

$ws.Range("B23").Value = 'NEARProtocol'
$ws.Range("C23").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D23").Value = '4.53'
$ws.Range("E23").Value = '  -2.77%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '71.23'
$ws.Range("E24").Value = '  +1.59%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("B26").Value = 'SuiNetwork'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D26").Value = '1.91'
$ws.Range("E26").Value = '  -5.32%  '

$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").Value = '9.91'
$ws.Range("E27").Value = '  -1.62%  '

$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '2.622.74'
$ws.Range("E29").Value = '  -1.88%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0956'
$ws.Range("E30").Value = '  -3.47%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '8.36'
$ws.Range("E31").Value = '  +2.21%  '

$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '529.48'
$ws.Range("E32").Value = '  -0.71%  '

$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '1.31'
$ws.Range("E33").Value = '  -3.70%  '

$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").Value = '1.85'
$ws.Range("E34").Value = '  +0.14%  '

$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.127'
$ws.Range("E36").Value = '  -4.02%  '

$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '157.97'
$ws.Range("E37").Value = '  +0.12%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '1.42'
$ws.Range("E38").Value = '  -3.73%  '

$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").Value = '18.66'
$ws.Range("E39").Value = '  -0.39%  '

$ws.Range("B40").Value = 'WhiteBITCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D40").Value = '18.62'
$ws.Range("E40").Value = '  +0.85%  '

$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D41").Value = '0.348'
$ws.Range("E41").Value = '  -2.62%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '1.77'
$ws.Range("E42").Value = '  -1.55%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").Value = '5.09'
$ws.Range("E43").Value = '  -1.53%  '

$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.29%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '2.48'
$ws.Range("E45").Value = '  -1.81%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '144.42'
$ws.Range("E46").Value = '  -3.90%  '

$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0275'
$ws.Range("E47").Value = '  -1.70%  '

$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").Value = '3.66'
$ws.Range("E48").Value = '  -1.32%  '

$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = '0.543'
$ws.Range("E49").Value = '  -3.48%  '

$ws.Range("B50").Value = 'Optimism'
$ws.Range("C50").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D50").Value = '1.68'
$ws.Range("E50").Value = '  -1.77%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.0749'
$ws.Range("E51").Value = '  -1.73%  '
